$d = $word.ActiveDocument

# Update the date line at the top of the document.
$d.Content.Find.Execute("2025-04-25 Friday", $false, $false, $false, $false, $false, $true, 1, $false, "2025-04-26 Saturday", 2)

# Update the division problems in the table, cell by cell (row, column)
# so that duplicate problem texts (e.g. "92÷9=" appears twice) are each
# replaced with their own distinct target value.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "87÷8="
$t.Cell(1, 2).Range.Text = "94÷3="
$t.Cell(1, 3).Range.Text = "98÷2="
$t.Cell(1, 4).Range.Text = "29÷9="
$t.Cell(1, 5).Range.Text = "58÷4="

$t.Cell(5, 1).Range.Text = "18÷4="
$t.Cell(5, 2).Range.Text = "69÷8="
$t.Cell(5, 3).Range.Text = "22÷6="
$t.Cell(5, 4).Range.Text = "24÷8="
$t.Cell(5, 5).Range.Text = "28÷4="

$t.Cell(9, 1).Range.Text = "75÷2="
$t.Cell(9, 2).Range.Text = "69÷4="
$t.Cell(9, 3).Range.Text = "56÷5="
$t.Cell(9, 4).Range.Text = "18÷3="
$t.Cell(9, 5).Range.Text = "18÷4="

$t.Cell(13, 1).Range.Text = "60÷3="
$t.Cell(13, 2).Range.Text = "49÷9="
$t.Cell(13, 3).Range.Text = "22÷4="
$t.Cell(13, 4).Range.Text = "39÷7="
$t.Cell(13, 5).Range.Text = "87÷2="

$t.Cell(17, 1).Range.Text = "79÷2="
$t.Cell(17, 2).Range.Text = "28÷2="
$t.Cell(17, 3).Range.Text = "37÷3="
$t.Cell(17, 4).Range.Text = "27÷7="
$t.Cell(17, 5).Range.Text = "72÷9="
